# "updated trading sheet and other models"
#
# On the "Main" worksheet (the sheet tab that is actually selected/active in
# the file -- it is named "Main" even though it holds the per-company
# trading comps table), a new comp row for JD Health (6618 HK) is inserted
# right after the Amerisource row (row 17), pushing every row below it down
# by one. The Amerisource row (17) also picks up the column-A flag value
# ("x") that all the other comp rows above it already had, and the new
# JD Health row gets the same flag.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Insert a new blank row above the old row 18 (LabCorp). This shifts rows
# 18-27 down to 19-28 and auto-adjusts the formulas/relative refs that live
# in the last row (old row 27, the Syapse row, which becomes row 28).
$ws.Rows.Item(18).Insert()

# Row 17 (Amerisource) previously had no value in column A; give it the same
# "x" flag as every other row in the table.
$ws.Range("A17").Value = "x"

# Fill in the newly inserted row 18 with the JD Health comp.
$ws.Range("A18").Value = "x"
$ws.Range("B18").Value = "JD Health"
$ws.Range("C18").Value = "6618 HK"

# The row insert shifts cell contents/formulas automatically, but this
# engine's Hyperlinks collection is not shift-aware, and clearing a single
# hyperlink only works by clearing the whole per-sheet collection -- so
# rebuild all three external hyperlinks, pointing the one that used to sit
# on B27 (Syapse -> SNCE.xlsx) at its new home, B28.
$ws.Range("B3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B28"), "SNCE.xlsx")
$ws.Hyperlinks.Add($ws.Range("B3"), "UNH.xlsx")
$ws.Hyperlinks.Add($ws.Range("B4"), "CVS.xlsx")

# Re-adding the hyperlinks restyles the cells; put them back on the shared
# "Hyperlink" cell style so B3/B4/B28 keep looking like they did before.
$ws.Range("B28").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"

# Match the saved selection/active cell from the edit.
$ws.Range("A19").Select()
